# Updates cryptos list with refreshed prices / volume percentages.
# Also fixes a row-order swap between EnergySwap (row 45) and PancakeSwap (row 46).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'20.554.70"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +1.63%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'1.475.24"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +2.57%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = "'1.009"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.27%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'0.9511"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +4.16%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'278.57"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +0.29%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'0.3623"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -1.25%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'0.3058"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -2.33%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'39.40"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +1.41%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'1.058"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +3.89%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'0.06654"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +1.86%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'1.003"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -0.21%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'5.526"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +2.28%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'18.11"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +3.25%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'6.197"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +2.14%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'0.9510"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +3.58%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'0.00001030"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +1.05%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'1.473.78"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +2.06%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'0.05941"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +5.92%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'69.37"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +3.19%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'5.503"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +1.81%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'14.49"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +0.53%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'11.13"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +1.83%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'2.261"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +0.15%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'20.592.73"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +1.51%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'142.85"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +5.70%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'2.125"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -3.28%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = "'17.19"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +1.82%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = "'1.632.83"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +2.36%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = "'113.73"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +3.04%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'3.950"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +6.44%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Value = "'5.033"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +3.24%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = "'0.8110"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -1.01%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = "'0.07979"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +4.67%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Value = "'1.513"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +2.54%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("E36").Value = "'  +6.46%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = "'0.05865"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -1.65%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").Value = "'4.736"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +0.82%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "'0.02050"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +2.64%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'10.39"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +1.28%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'0.9518"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +3.06%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'0.1880"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +2.84%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'7.422"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +8.46%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'0.5301"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +1.30%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("B45").Value = "'PancakeSwap"
$ws.Range("B45").Style = "Normal"
$ws.Range("C45").Value = "'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("C45").Style = "Normal"
$ws.Range("D45").Value = "'3.533"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +0.05%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("B46").Value = "'EnergySwap"
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").Value = "'12.25"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +2.81%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'118.32"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -1.00%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'0.5198"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +1.07%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'1.817"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +3.08%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'0.06479"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +2.27%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'0.9789"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -1.47%  "
$ws.Range("E51").Style = "Normal"
